$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing existing rows 14-27 down to 15-28.
# Excel's default Insert() copies formatting from the row above, which is
# what we want (column D keeps its date number-format style).
$ws.Rows(14).Insert()

# Populate the new row 14 with the weekly data point (same categorical
# fields as the old row that is now row 15, but a new date + prices).
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = 44771
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 100112043
$ws.Range("G14").Value = "Pepino dulce"
$ws.Range("H14").Value = "Cultivar XV región"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 9000
$ws.Range("M14").Value = 8500
$ws.Range("N14").Value = "$/caja 10 kilos"
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 850
$ws.Range("Q14").Value = 10
$ws.Range("R14").Value = "Hortaliza"
